$wb = $excel.ActiveWorkbook

# ALC row 40: Stuck in the Moment | Horn Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1230.9286
$ws.Range("I40").Value = 1081.8182
$ws.Range("K40").Value = 1081.8182
$ws.Range("M40").Value = -906.8181999999999

# ALC row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 11832.9
$ws.Range("I86").Value = 2068.6
$ws.Range("J86").Value = 21597.2
$ws.Range("K86").Value = 2068.6
$ws.Range("L86").Value = 21597.2
$ws.Range("M86").Value = -945.5999999999999
$ws.Range("N86").Value = -23843.2

# ALC row 87: There Was a Late Fee | Noble Gold
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 38687.332
$ws.Range("J87").Value = 38687.332
$ws.Range("L87").Value = 38687.332
$ws.Range("N87").Value = -41183.332

# ALC row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 11832.9
$ws.Range("I89").Value = 2068.6
$ws.Range("J89").Value = 21597.2
$ws.Range("K89").Value = 10343
$ws.Range("L89").Value = 107986
$ws.Range("M89").Value = -4727
$ws.Range("N89").Value = -119218

# ALC row 90: A Gate Arcane Is Dragon's Bane (L) | Noble Gold
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 38687.332
$ws.Range("J90").Value = 38687.332
$ws.Range("L90").Value = 116061.996
$ws.Range("N90").Value = -128541.996

# ALC row 127: Liquid Competence | Competent Craftsman's Draught
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1776
$ws.Range("I127").Value = 897
$ws.Range("J127").Value = 2069
$ws.Range("K127").Value = 2691
$ws.Range("L127").Value = 6207
$ws.Range("M127").Value = 2269
$ws.Range("N127").Value = -16127

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3469
$ws.Range("I132").Value = 3532.6086
$ws.Range("K132").Value = 10597.8258
$ws.Range("M132").Value = -8067.825800000001

# ARM row 2: Ain't Got No Ingots | Bronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1550.3438
$ws.Range("I2").Value = 993.5217
$ws.Range("K2").Value = 993.5217
$ws.Range("M2").Value = -880.5217

# ARM row 53: Metal Fatigue | Mythril Vambraces
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 20000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 20000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 20000
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -21364

# ARM row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1504246.9
$ws.Range("I61").Value = 3003185.2
$ws.Range("J61").Value = 5308.3335
$ws.Range("K61").Value = 3003185.2
$ws.Range("L61").Value = 5308.3335
$ws.Range("M61").Value = -3002973.2
$ws.Range("N61").Value = -5732.3335

# ARM row 63: Rivets Run through It | Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2485.7144
$ws.Range("I63").Value = 2485.7144
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2485.7144
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1799.7144
$ws.Range("N63").ClearContents()

# ARM row 66: A Riveting Revival (L) | Mythrite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2485.7144
$ws.Range("I66").Value = 2485.7144
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12428.572
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -8996.572
$ws.Range("N66").ClearContents()

# ARM row 116: No Scope | Titanbronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1550.3438
$ws.Range("I116").Value = 993.5217
$ws.Range("K116").Value = 993.5217
$ws.Range("M116").Value = 1300.4783

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 17928.469
$ws.Range("I132").Value = 2002.375
$ws.Range("J132").Value = 33854.562
$ws.Range("K132").Value = 6007.125
$ws.Range("L132").Value = 101563.686
$ws.Range("M132").Value = -3477.125
$ws.Range("N132").Value = -106623.686

# ARM row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1504246.9
$ws.Range("I136").Value = 3003185.2
$ws.Range("J136").Value = 5308.3335
$ws.Range("K136").Value = 9009555.600000001
$ws.Range("L136").Value = 15925.0005
$ws.Range("M136").Value = -9007005.600000001
$ws.Range("N136").Value = -21025.0005

# BSM row 3: Hells Bells | Bronze Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1550.3438
$ws.Range("I3").Value = 993.5217
$ws.Range("K3").Value = 993.5217
$ws.Range("M3").Value = -879.5217

# BSM row 82: Spirituality Inspector | Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 22646
$ws.Range("I82").Value = 6381.5713
$ws.Range("J82").Value = 51108.75
$ws.Range("K82").Value = 6381.5713
$ws.Range("L82").Value = 51108.75
$ws.Range("M82").Value = -5998.5713
$ws.Range("N82").Value = -51874.75

# BSM row 85: The Clamor for Hammers (L) | Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 22646
$ws.Range("I85").Value = 6381.5713
$ws.Range("J85").Value = 51108.75
$ws.Range("K85").Value = 6381.5713
$ws.Range("L85").Value = 51108.75
$ws.Range("M85").Value = -5055.5713
$ws.Range("N85").Value = -53760.75

# BSM row 107: The Gold Experience | Deepgold Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4262.6
$ws.Range("I107").Value = 2150
$ws.Range("J107").Value = 5671
$ws.Range("K107").Value = 2150
$ws.Range("L107").Value = 5671
$ws.Range("M107").Value = -230
$ws.Range("N107").Value = -9511

# BSM row 134: Ruthenium Supremium | Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3868.238
$ws.Range("I134").Value = 4128.0527
$ws.Range("K134").Value = 12384.1581
$ws.Range("M134").Value = -9849.158100000001

# CRP row 22: Driving Up the Wall | Elm Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1320
$ws.Range("I22").Value = 933.3333
$ws.Range("K22").Value = 933.3333
$ws.Range("M22").Value = -583.3333

# CRP row 31: Wall Not Found | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8051
$ws.Range("I31").Value = 8944.975
$ws.Range("J31").Value = 3581.125
$ws.Range("K31").Value = 8944.975
$ws.Range("L31").Value = 3581.125
$ws.Range("M31").Value = -8649.975
$ws.Range("N31").Value = -4171.125

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8051
$ws.Range("I34").Value = 8944.975
$ws.Range("J34").Value = 3581.125
$ws.Range("K34").Value = 8944.975
$ws.Range("L34").Value = 3581.125
$ws.Range("M34").Value = -8742.975
$ws.Range("N34").Value = -3985.125

# CRP row 55: Ready for a Rematch | Mythril Lance
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 16666.666
$ws.Range("J55").Value = 16666.666
$ws.Range("L55").Value = 16666.666
$ws.Range("N55").Value = -17296.666

# CRP row 105: Zelkova, My Love | Zelkova Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 12501271
$ws.Range("I105").Value = 20833780
$ws.Range("J105").Value = 2508.25
$ws.Range("K105").Value = 20833780
$ws.Range("L105").Value = 2508.25
$ws.Range("M105").Value = -20832033
$ws.Range("N105").Value = -6002.25

# CUL row 37: I Love Lamprey | Eel Pie
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 19293616
$ws.Range("J37").Value = 19293616
$ws.Range("L37").Value = 57880848
$ws.Range("N37").Value = -57881072

# CUL row 131: The Mountain Steeped | Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 741.58
$ws.Range("I131").Value = 310
$ws.Range("J131").Value = 754.92786
$ws.Range("K131").Value = 930
$ws.Range("L131").Value = 2264.78358
$ws.Range("M131").Value = 4110
$ws.Range("N131").Value = -12344.78358

# CUL row 136: Simple Is Hardest | Spaghetti al Olio e Peperoncino
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 4930
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# GSM row 113: Copious Crystal Cannons | Manasilver Nugget
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3487.375
$ws.Range("J113").Value = 4525
$ws.Range("L113").Value = 4525
$ws.Range("N113").Value = -8865

# GSM row 126: Gold Rush Order | Phrygian Gold Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4003.2285
$ws.Range("I126").Value = 2931.5264
$ws.Range("J126").Value = 5275.875
$ws.Range("K126").Value = 8794.5792
$ws.Range("L126").Value = 15827.625
$ws.Range("M126").Value = -6324.5792
$ws.Range("N126").Value = -20767.625

# GSM row 132: On Board for Lar | Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 33758.41
$ws.Range("I132").Value = 4407.9165
$ws.Range("K132").Value = 13223.7495
$ws.Range("M132").Value = -10693.7495

# LTW row 9: From the Sands to the Stage | Leather Himantes
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 355
$ws.Range("I9").Value = 355
$ws.Range("K9").Value = 355
$ws.Range("M9").Value = -131

# LTW row 40: Best Served Toad | Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4849.6665
$ws.Range("I40").Value = 4021.7778
$ws.Range("K40").Value = 4021.7778
$ws.Range("M40").Value = -3885.7778

# LTW row 46: Supply Side Logic | Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2515
$ws.Range("I46").Value = 2472.5
$ws.Range("J46").Value = 2600
$ws.Range("K46").Value = 2472.5
$ws.Range("L46").Value = 2600
$ws.Range("M46").Value = -2284.5
$ws.Range("N46").Value = -2976

# LTW row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 170.05556
$ws.Range("I55").Value = 122.166664
$ws.Range("J55").Value = 194
$ws.Range("K55").Value = 122.166664
$ws.Range("L55").Value = 194
$ws.Range("M55").Value = 50.833336
$ws.Range("N55").Value = -540

# LTW row 61: Spelling Me Softly | Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4497.476
$ws.Range("I61").Value = 1814.7693
$ws.Range("J61").Value = 8856.875
$ws.Range("K61").Value = 1814.7693
$ws.Range("L61").Value = 8856.875
$ws.Range("M61").Value = -1612.7693
$ws.Range("N61").Value = -9260.875

# LTW row 113: Peace in Rest | Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4497.476
$ws.Range("I113").Value = 1814.7693
$ws.Range("J113").Value = 8856.875
$ws.Range("K113").Value = 1814.7693
$ws.Range("L113").Value = 8856.875
$ws.Range("M113").Value = 355.2307000000001
$ws.Range("N113").Value = -13196.875

# WVR row 122: Heavy Armoire | Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1666.6818
$ws.Range("I122").Value = 1572
$ws.Range("J122").Value = 1988.6
$ws.Range("K122").Value = 4716
$ws.Range("L122").Value = 5965.799999999999
$ws.Range("M122").Value = -2266
$ws.Range("N122").Value = -10865.8

# WVR row 126: A Polished Purchase | Snow Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1490.1428
$ws.Range("I126").Value = 951.3333
$ws.Range("J126").Value = 2460
$ws.Range("K126").Value = 2853.9999
$ws.Range("L126").Value = 7380
$ws.Range("M126").Value = -383.9998999999998
$ws.Range("N126").Value = -12320
